$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.764.64'
$ws.Range('E2').Value = '  -1.97%  '

$ws.Range('D3').Value = '3.426.04'
$ws.Range('E3').Value = '  -3.07%  '

$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.08%  '

$ws.Range('D5').Value = '581.24'
$ws.Range('E5').Value = '  -3.72%  '

$ws.Range('D6').Value = '135.03'
$ws.Range('E6').Value = '  -6.26%  '

$ws.Range('E7').Value = '  +0.08%  '

$ws.Range('D8').Value = '3.425.41'
$ws.Range('E8').Value = '  -3.23%  '

$ws.Range('D9').Value = '0.483'
$ws.Range('E9').Value = '  -5.81%  '

$ws.Range('D10').Value = '0.120'
$ws.Range('E10').Value = '  -7.96%  '

$ws.Range('D11').Value = '7.06'
$ws.Range('E11').Value = '  -9.66%  '

$ws.Range('D12').Value = '0.375'
$ws.Range('E12').Value = '  -7.84%  '

$ws.Range('D13').Value = '4.000.85'
$ws.Range('E13').Value = '  -2.98%  '

$ws.Range('D14').Value = '0.0000179'
$ws.Range('E14').Value = '  -7.88%  '

$ws.Range('D15').Value = '3.423.84'
$ws.Range('E15').Value = '  -3.15%  '

$ws.Range('E16').Value = '  -1.65%  '

$ws.Range('B17').Value = 'Avalanche'
$ws.Range('C17').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D17').Value = '26.15'
$ws.Range('E17').Value = '  -7.74%  '

$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '64.630.25'
$ws.Range('E18').Value = '  -1.89%  '

$ws.Range('D19').Value = '9.47'
$ws.Range('E19').Value = '  -13.42%  '

$ws.Range('D20').Value = '5.75'
$ws.Range('E20').Value = '  -7.16%  '

$ws.Range('D21').Value = '13.49'
$ws.Range('E21').Value = '  -7.11%  '

$ws.Range('D22').Value = '380.84'
$ws.Range('E22').Value = '  -9.35%  '

$ws.Range('B23').Value = 'Polygon'
$ws.Range('C23').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D23').Value = '0.542'
$ws.Range('E23').Value = '  -8.81%  '

$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D24').Value = '1.00'
$ws.Range('E24').Value = '  -0.04%  '

$ws.Range('D25').Value = '71.62'
$ws.Range('E25').Value = '  -6.76%  '

$ws.Range('D26').Value = '3.554.56'
$ws.Range('E26').Value = '  -3.06%  '

$ws.Range('D27').Value = '0.0000105'
$ws.Range('E27').Value = '  -8.46%  '

$ws.Range('D28').Value = '0.994'
$ws.Range('E28').Value = '  -0.67%  '

$ws.Range('D29').Value = '7.22'
$ws.Range('E29').Value = '  -7.27%  '

$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').Value = '7.98'
$ws.Range('E30').Value = '  -10.10%  '

$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '2.18'
$ws.Range('E31').Value = '  -11.54%  '

$ws.Range('D32').Value = '3.437.95'
$ws.Range('E32').Value = '  -2.66%  '

$ws.Range('E33').Value = '  +0.04%  '

$ws.Range('D34').Value = '22.90'
$ws.Range('E34').Value = '  -5.44%  '

$ws.Range('D35').Value = '0.141'
$ws.Range('E35').Value = '  -9.11%  '

$ws.Range('D36').Value = '168.52'
$ws.Range('E36').Value = '  -4.79%  '

$ws.Range('D37').Value = '1.18'
$ws.Range('E37').Value = '  -11.96%  '

$ws.Range('D38').Value = '6.69'
$ws.Range('E38').Value = '  -11.62%  '

$ws.Range('D39').Value = '1.45'
$ws.Range('E39').Value = '  -10.69%  '

$ws.Range('D40').Value = '4.61'
$ws.Range('E40').Value = '  -11.74%  '

$ws.Range('D41').Value = '0.0758'
$ws.Range('E41').Value = '  -7.03%  '

$ws.Range('E42').Value = '  -6.56%  '

$ws.Range('D43').Value = '0.998'
$ws.Range('E43').Value = '  -0.06%  '

$ws.Range('D44').Value = '41.82'
$ws.Range('E44').Value = '  -8.11%  '

$ws.Range('D45').Value = '4.28'
$ws.Range('E45').Value = '  -13.85%  '

$ws.Range('D46').Value = '1.62'
$ws.Range('E46').Value = '  -8.72%  '

$ws.Range('D47').Value = '1.12'
$ws.Range('E47').Value = '  +0.88%  '

$ws.Range('D48').Value = '22.58'
$ws.Range('E48').Value = '  -3.74%  '

$ws.Range('D49').Value = '6.45'
$ws.Range('E49').Value = '  -8.42%  '

$ws.Range('D50').Value = '2.185.91'
$ws.Range('E50').Value = '  -5.87%  '

$ws.Range('D51').Value = '1.99'
$ws.Range('E51').Value = '  -15.97%  '
